# Insert two new data rows into the "Nectarín" sheet right after the existing
# row 396 (i.e. at what is currently row 397), pushing every following row
# down by two positions. Then populate the two new rows with the new
# "Artic Mist" price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 397; Excel will shift rows 397.. downward.
$ws.Rows.Item(397).Insert()
$ws.Rows.Item(397).Insert()

# --- New row 397 ---------------------------------------------------------
$ws.Range("A397").Value = 7
$ws.Range("B397").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C397").Value = "Ñuble"
$ws.Range("D397").Value = 44995
$ws.Range("E397").Value = 16
$ws.Range("F397").Value = "Fruta"
$ws.Range("G397").Value = 100103
$ws.Range("H397").Value = "Frutos de hueso (carozo)"
$ws.Range("I397").Value = 100103006
$ws.Range("J397").Value = "Nectarín"
$ws.Range("K397").Value = "Artic Mist"
$ws.Range("L397").Value = "Especial"
$ws.Range("M397").Value = 60
$ws.Range("N397").Value = 15000
$ws.Range("O397").Value = 15000
$ws.Range("P397").Value = 15000
$ws.Range("Q397").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R397").Value = "Región de O'Higgins"
$ws.Range("S397").Value = 1000
$ws.Range("T397").Value = 15

# --- New row 398 ---------------------------------------------------------
$ws.Range("A398").Value = 7
$ws.Range("B398").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C398").Value = "Ñuble"
$ws.Range("D398").Value = 44995
$ws.Range("E398").Value = 16
$ws.Range("F398").Value = "Fruta"
$ws.Range("G398").Value = 100103
$ws.Range("H398").Value = "Frutos de hueso (carozo)"
$ws.Range("I398").Value = 100103006
$ws.Range("J398").Value = "Nectarín"
$ws.Range("K398").Value = "Artic Mist"
$ws.Range("L398").Value = "Primera"
$ws.Range("M398").Value = 60
$ws.Range("N398").Value = 13000
$ws.Range("O398").Value = 13000
$ws.Range("P398").Value = 13000
$ws.Range("Q398").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R398").Value = "Región de O'Higgins"
$ws.Range("S398").Value = 867
$ws.Range("T398").Value = 15
